$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'37.123.64"
$ws.Cells.Item(2, 5).Value = "  +1.28%  "

$ws.Cells.Item(3, 4).Value = "'2.057.43"
$ws.Cells.Item(3, 5).Value = "  -2.57%  "

$ws.Cells.Item(4, 5).Value = "  +0.09%  "

$ws.Cells.Item(5, 4).Value = "'249.06"
$ws.Cells.Item(5, 5).Value = "  -1.79%  "

$ws.Cells.Item(6, 5).Value = "  -1.05%  "

$ws.Cells.Item(7, 5).Value = "  -0.04%  "

$ws.Cells.Item(8, 4).Value = "'55.63"
$ws.Cells.Item(8, 5).Value = "  +15.16%  "

$ws.Cells.Item(9, 4).Value = "'61.77"
$ws.Cells.Item(9, 5).Value = "  +3.81%  "

$ws.Cells.Item(10, 4).Value = "'0.380"
$ws.Cells.Item(10, 5).Value = "  +1.77%  "

$ws.Cells.Item(11, 4).Value = "'0.0794"
$ws.Cells.Item(11, 5).Value = "  +6.35%  "

$ws.Cells.Item(12, 5).Value = "  +5.72%  "

$ws.Cells.Item(13, 4).Value = "'15.19"
$ws.Cells.Item(13, 5).Value = "  +6.13%  "

$ws.Cells.Item(14, 4).Value = "'2.357.55"
$ws.Cells.Item(14, 5).Value = "  -2.59%  "

$ws.Cells.Item(15, 4).Value = "'0.818"
$ws.Cells.Item(15, 5).Value = "  -1.61%  "

$ws.Cells.Item(16, 4).Value = "'5.25"
$ws.Cells.Item(16, 5).Value = "  +2.65%  "

$ws.Cells.Item(17, 4).Value = "'2.059.45"
$ws.Cells.Item(17, 5).Value = "  -2.38%  "

$ws.Cells.Item(18, 4).Value = "'37.061.16"
$ws.Cells.Item(18, 5).Value = "  +1.13%  "

$ws.Cells.Item(19, 4).Value = "'0.0₃0918"
$ws.Cells.Item(19, 5).Value = "  +9.96%  "

$ws.Cells.Item(20, 4).Value = "'72.41"
$ws.Cells.Item(20, 5).Value = "  -1.49%  "

$ws.Cells.Item(21, 4).Value = "'14.24"
$ws.Cells.Item(21, 5).Value = "  +7.45%  "

$ws.Cells.Item(22, 4).Value = "'5.38"
$ws.Cells.Item(22, 5).Value = "  +3.55%  "

$ws.Cells.Item(23, 4).Value = "'237.08"

$ws.Cells.Item(24, 5).Value = "  -0.04%  "

$ws.Cells.Item(25, 4).Value = "'2.43"
$ws.Cells.Item(25, 5).Value = "  -2.02%  "

$ws.Cells.Item(26, 4).Value = "'170.08"
$ws.Cells.Item(26, 5).Value = "  -1.52%  "

$ws.Cells.Item(27, 5).Value = "  -1.46%  "

$ws.Cells.Item(28, 4).Value = "'20.26"
$ws.Cells.Item(28, 5).Value = "  -6.06%  "

$ws.Cells.Item(29, 5).Value = "  -2.03%  "

$ws.Cells.Item(30, 5).Value = "  -0.35%  "

$ws.Cells.Item(31, 4).Value = "'4.54"
$ws.Cells.Item(31, 5).Value = "  +1.22%  "

$ws.Cells.Item(32, 5).Value = "  +9.95%  "

$ws.Cells.Item(33, 4).Value = "'0.0623"
$ws.Cells.Item(33, 5).Value = "  +3.21%  "

$ws.Cells.Item(34, 4).Value = "'4.33"
$ws.Cells.Item(34, 5).Value = "  +5.62%  "

$ws.Cells.Item(35, 5).Value = "  +0.14%  "

$ws.Cells.Item(36, 4).Value = "'0.0862"
$ws.Cells.Item(36, 5).Value = "  -5.62%  "

$ws.Cells.Item(37, 5).Value = "  -3.63%  "

$ws.Cells.Item(38, 4).Value = "'1.76"
$ws.Cells.Item(38, 5).Value = "  -6.69%  "

$ws.Cells.Item(39, 5).Value = "  +1.53%  "

$ws.Cells.Item(40, 5).Value = "  +22.49%  "

$ws.Cells.Item(41, 4).Value = "'17.95"
$ws.Cells.Item(41, 5).Value = "  +11.98%  "

$ws.Cells.Item(42, 5).Value = "  -0.84%  "

$ws.Cells.Item(43, 5).Value = "  -3.58%  "

$ws.Cells.Item(44, 4).Value = "'96.32"
$ws.Cells.Item(44, 5).Value = "  -2.37%  "

$ws.Cells.Item(45, 4).Value = "'4.34"
$ws.Cells.Item(45, 5).Value = "  +46.14%  "

$ws.Cells.Item(46, 5).Value = "  +0.50%  "

$ws.Cells.Item(47, 4).Value = "'14.08"
$ws.Cells.Item(47, 5).Value = "  -52.67%  "

$ws.Cells.Item(48, 5).Value = "  +5.97%  "

$ws.Cells.Item(49, 4).Value = "'1.297.98"
$ws.Cells.Item(49, 5).Value = "  -3.55%  "

$ws.Cells.Item(50, 4).Value = "'2.90"
$ws.Cells.Item(50, 5).Value = "  +2.46%  "

$ws.Cells.Item(51, 2).Value = "FraxShare"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(51, 4).Value = "'6.79"
$ws.Cells.Item(51, 5).Value = "  -5.62%  "
